$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GL")

$xlPasteFormats = -4122

# ---- Copy existing cell formatting onto the new cells first ----
# Dates with "m/d/yyyy" style (same as A56/A59/A64 pattern)
$ws.Range("A59").Copy()
$ws.Range("A60").PasteSpecial($xlPasteFormats)
$ws.Range("A63").PasteSpecial($xlPasteFormats)
$ws.Range("A64").PasteSpecial($xlPasteFormats)
$ws.Range("A65").PasteSpecial($xlPasteFormats)
$ws.Range("A66").PasteSpecial($xlPasteFormats)

# Dates with "d-mmm" style (same as A49 pattern)
$ws.Range("A49").Copy()
$ws.Range("A61:A62").PasteSpecial($xlPasteFormats)

# Amount / rate columns (style used throughout rows 49-59)
$ws.Range("B49").Copy()
$ws.Range("B60:B66").PasteSpecial($xlPasteFormats)

$ws.Range("C49").Copy()
$ws.Range("C60:C66").PasteSpecial($xlPasteFormats)

# H56 gets the same currency style used elsewhere in column H (e.g. H24)
$ws.Range("H24").Copy()
$ws.Range("H56").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---- Populate values ----
# Row 60
$ws.Range("A60").Value = 42886
$ws.Range("B60").Value = 2500
$ws.Range("C60").Value = 138.21
$ws.Range("D60").Value = "Quincena "

# Row 61
$ws.Range("A61").Value = 42894
$ws.Range("B61").Value = 2300
$ws.Range("C61").Value = 128.97
$ws.Range("D61").Value = "`$800 pintura, `$700 fumigada y `$800 regalo medicina Meli"

# Row 62
$ws.Range("A62").Value = 42894
$ws.Range("B62").Value = -800
$ws.Range("C62").Value = -44.444444444444443
$ws.Range("D62").Value = "Regalo medicina Meli"

# Row 63
$ws.Range("A63").Value = 42900
$ws.Range("B63").Value = 7800
$ws.Range("C63").Value = 438.39
$ws.Range("D63").Value = "`$2,500 quincena, `$4,500 y `$800 para arreglar depa?"

# Row 64
$ws.Range("A64").Value = 42913
$ws.Range("B64").Value = 1000
$ws.Range("C64").Value = 57.14
$ws.Range("D64").Value = "Prestamo especial. Dr? Medicinas?"

# Row 65
$ws.Range("A65").Value = 42916
$ws.Range("B65").Value = 2500
$ws.Range("C65").Value = 141.87
$ws.Range("D65").Value = "Quincena"

# Row 66
$ws.Range("A66").Value = 42921
$ws.Range("B66").Value = 15000
$ws.Range("C66").Value = 833.91
$ws.Range("D66").Value = "`$9,300 predial, saldo y multa de agua y posiblemente mantenimient (`$4,500)"

# C59 loses its formula but keeps the same cached value
$ws.Range("C59").Value = -16.286644951140065

# Update the remembered selection on the active sheet
$ws.Range("G61").Select()

$wb.Save()
